$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.442.56'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.438.37'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.31%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '407.79'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.43'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.48%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.45%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.127'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.17'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.00%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.91%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.50'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.00'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.63%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.477.46'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '62.452.19'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '11.45'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.83%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000140'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.18'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -5.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '83.86'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '314.14'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.92'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.17'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.77'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +8.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '29.78'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.19'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.07%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +5.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.64'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.52%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.08%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '42.94'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.61%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.42'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0487'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '51.47'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.328'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +13.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.95'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '138.34'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.52%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.35%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.86'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.23'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '21.57'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.127.66'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.60%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.33%  '
$ws.Range('B50').Value = 'Fetch.AI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.76'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +21.12%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.94'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.62%  '
